$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cooked-food recipes added, first pass (ascending XP order, matches final row order)
$ws.Range("A5").Value = "Bread"
$ws.Range("A5").Font.Size = 11
$ws.Range("B5").Value = 13
$ws.Range("B5").Font.Size = 11
$ws.Range("C5").Value = 5
$ws.Range("C5").Font.Size = 11
$ws.Range("D5").Value = 8

$ws.Range("A6").Value = "Herring"
$ws.Range("A6").Font.Size = 11
$ws.Range("B6").Value = 15
$ws.Range("B6").Font.Size = 11
$ws.Range("C6").Value = 10
$ws.Range("C6").Font.Size = 11
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "Seahorse"
$ws.Range("A7").Font.Size = 11
$ws.Range("B7").Value = 30
$ws.Range("B7").Font.Size = 11
$ws.Range("C7").Value = 65
$ws.Range("C7").Font.Size = 11
$ws.Range("D7").Value = 4

$ws.Range("A8").Value = "Trout"
$ws.Range("A8").Font.Size = 11
$ws.Range("B8").Value = 25
$ws.Range("B8").Font.Size = 11
$ws.Range("C8").Value = 27
$ws.Range("C8").Font.Size = 11
$ws.Range("D8").Value = 4

$ws.Range("A9").Value = "Salmon"
$ws.Range("A9").Font.Size = 11
$ws.Range("B9").Value = 30
$ws.Range("B9").Font.Size = 11
$ws.Range("C9").Value = 58
$ws.Range("C9").Font.Size = 11
$ws.Range("D9").Value = 4

$ws.Range("A10").Value = "Lobster"
$ws.Range("A10").Font.Size = 11
$ws.Range("B10").Value = 40
$ws.Range("B10").Font.Size = 11
$ws.Range("C10").Value = 108
$ws.Range("C10").Font.Size = 11
$ws.Range("D10").Value = 4

$ws.Range("A11").Value = "Chicken"
$ws.Range("A11").Font.Size = 11
$ws.Range("B11").Value = 44
$ws.Range("B11").Font.Size = 11
$ws.Range("C11").Value = 8
$ws.Range("C11").Font.Size = 11
$ws.Range("D11").Value = 8

$ws.Range("A12").Value = "Swordfish"
$ws.Range("A12").Font.Size = 11
$ws.Range("B12").Value = 50
$ws.Range("B12").Font.Size = 11
$ws.Range("C12").Value = 134
$ws.Range("C12").Font.Size = 11
$ws.Range("D12").Value = 5

$ws.Range("A13").Value = "Anglerfish"
$ws.Range("A13").Font.Size = 11
$ws.Range("B13").Value = 60
$ws.Range("B13").Font.Size = 11
$ws.Range("C13").Value = 209
$ws.Range("C13").Font.Size = 11
$ws.Range("D13").Value = 6

$ws.Range("A14").Value = "Fanfish"
$ws.Range("A14").Font.Size = 11
$ws.Range("B14").Value = 60
$ws.Range("B14").Font.Size = 11
$ws.Range("C14").Value = 250
$ws.Range("C14").Font.Size = 11
$ws.Range("D14").Value = 6

$ws.Range("A15").Value = "Crab"
$ws.Range("A15").Font.Size = 11
$ws.Range("B15").Value = 70
$ws.Range("B15").Font.Size = 11
$ws.Range("C15").Value = 280
$ws.Range("C15").Font.Size = 11
$ws.Range("D15").Value = 7

$ws.Range("A16").Value = "Plain Pizza Slice"
$ws.Range("A16").Font.Size = 11
$ws.Range("B16").Value = 72
$ws.Range("B16").Font.Size = 11
$ws.Range("C16").Value = 2
$ws.Range("C16").Font.Size = 11
$ws.Range("D16").Value = 8

$ws.Range("A17").Value = "Carp"
$ws.Range("A17").Font.Size = 11
$ws.Range("B17").Value = 75
$ws.Range("B17").Font.Size = 11
$ws.Range("C17").Value = 395
$ws.Range("C17").Font.Size = 11
$ws.Range("D17").Value = 7

$ws.Range("A18").Value = "Shark"
$ws.Range("A18").Font.Size = 11
$ws.Range("B18").Value = 80
$ws.Range("B18").Font.Size = 11
$ws.Range("C18").Value = 674
$ws.Range("C18").Font.Size = 11
$ws.Range("D18").Value = 8

$ws.Range("A19").Value = "Cave Fish"
$ws.Range("A19").Font.Size = 11
$ws.Range("B19").Value = 100
$ws.Range("B19").Font.Size = 11
$ws.Range("C19").Value = 538
$ws.Range("C19").Font.Size = 11
$ws.Range("D19").Value = 9

$ws.Range("A20").Value = "Beef Pie"
$ws.Range("A20").Font.Size = 11
$ws.Range("B20").Value = 161
$ws.Range("B20").Font.Size = 11
$ws.Range("C20").Value = 37
$ws.Range("C20").Font.Size = 11
$ws.Range("D20").Value = 8

$ws.Range("A21").Value = "Manta Ray"
$ws.Range("A21").Font.Size = 11
$ws.Range("B21").Value = 125
$ws.Range("B21").Font.Size = 11
$ws.Range("C21").Value = 1624
$ws.Range("C21").Font.Size = 11
$ws.Range("D21").Value = 10

$ws.Range("A23").Value = "Whale"
$ws.Range("A23").Font.Size = 11
$ws.Range("B23").Value = 150
$ws.Range("B23").Font.Size = 11
$ws.Range("C23").Value = 2048
$ws.Range("C23").Font.Size = 11
$ws.Range("D23").Value = 11

$ws.Range("A24").Value = "Meat Pizza Slice"
$ws.Range("A24").Font.Size = 11
$ws.Range("B24").Value = 161
$ws.Range("B24").Font.Size = 11
$ws.Range("C24").Value = 4
$ws.Range("C24").Font.Size = 11
$ws.Range("D24").Value = 8

$ws.Range("A25").Value = "Strawberry Cupcake"
$ws.Range("A25").Font.Size = 11
$ws.Range("B25").Value = 214
$ws.Range("B25").Font.Size = 11
$ws.Range("C25").Value = 42
$ws.Range("C25").Font.Size = 11
$ws.Range("D25").Value = 8

$ws.Range("A27").Value = "Cherry Cupcake"
$ws.Range("A27").Font.Size = 11
$ws.Range("B27").Value = 252
$ws.Range("B27").Font.Size = 11
$ws.Range("C27").Value = 78
$ws.Range("C27").Font.Size = 11
$ws.Range("D27").Value = 8

$ws.Range("A29").Value = "Apple Pie"
$ws.Range("A29").Font.Size = 11
$ws.Range("B29").Value = 424
$ws.Range("B29").Font.Size = 11
$ws.Range("C29").Value = 253
$ws.Range("C29").Font.Size = 11
$ws.Range("D29").Value = 8

$ws.Range("A30").Value = "Strawberry Cake"
$ws.Range("A30").Font.Size = 11
$ws.Range("B30").Value = 490
$ws.Range("B30").Font.Size = 11
$ws.Range("C30").Value = 378
$ws.Range("C30").Font.Size = 11
$ws.Range("D30").Value = 8

$ws.Range("A31").Value = "Carrot Cake"
$ws.Range("A31").Font.Size = 11
$ws.Range("B31").Value = 603
$ws.Range("B31").Font.Size = 11
$ws.Range("C31").Value = 751
$ws.Range("C31").Font.Size = 11
$ws.Range("D31").Value = 8

# Soup recipes added in a second pass (ascending XP order, matches final row order)
$ws.Range("A22").Value = "Basic Soup"
$ws.Range("A22").Font.Size = 11
$ws.Range("B22").Value = 126
$ws.Range("B22").Font.Size = 11
$ws.Range("C22").Value = 8
$ws.Range("C22").Font.Size = 11
$ws.Range("D22").Value = 7

$ws.Range("A26").Value = "Hearty Soup"
$ws.Range("A26").Font.Size = 11
$ws.Range("B26").Value = 223
$ws.Range("B26").Font.Size = 11
$ws.Range("C26").Value = 22
$ws.Range("C26").Font.Size = 11
$ws.Range("D26").Value = 7

$ws.Range("A28").Value = "Cream Corn Soup"
$ws.Range("A28").Font.Size = 11
$ws.Range("B28").Value = 274
$ws.Range("B28").Font.Size = 11
$ws.Range("C28").Value = 42
$ws.Range("C28").Font.Size = 11
$ws.Range("D28").Value = 7

$ws.Range("A32").Value = "Chicken Soup"
$ws.Range("A32").Font.Size = 11
$ws.Range("B32").Value = 617
$ws.Range("B32").Font.Size = 11
$ws.Range("C32").Value = 96
$ws.Range("C32").Font.Size = 11
$ws.Range("D32").Value = 7

# Column A widened to fit the new, longer recipe names
$ws.Columns("A:A").ColumnWidth = 15.5

# Restore the selection to the cell that was active when the edit finished
$ws.Range("C11").Select()
